# Auto-generated Excel COM-interop edit script
# Applies the cell-value updates described by the diff to the active worksheet.
# (cryptos.xlsx - 'cryptocurrency ranking' snapshot refresh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.024.67"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "2.663.51"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'525.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'144.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").Value = "'6.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.02%  "
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").Value = "'0.336"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "3.131.92"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "59.008.98"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "'21.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "2.671.41"
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "'338.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'4.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").Value = "'10.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "'6.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'64.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("D24").Value = "'0.420"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").Value = "'0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").Value = "0.0₃0802"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "'6.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'18.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "'150.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("D36").Value = "'0.889"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.32%  "
$ws.Range("D37").Value = "'0.872"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'36.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  -6.06%  "
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'19.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "'275.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").Value = "'0.0970"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0532"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.046.65"
$ws.Range("E48").Value = "  -3.84%  "
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0228"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'18.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.37%  "
